$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "*norm_c" column (old F),
# shifting it and everything after it two columns to the right.
$ws.Columns("F:G").Insert()

# Give the two new columns the same width as their neighbouring data
# columns (D:E), matching the author's resize after inserting them.
$ws.Columns("F:G").ColumnWidth = 13

# New column F: "lepton beam" header + "e_minus" for every data row.
$ws.Range("F1").Value = "lepton beam"
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 6).Value = "e_minus"
}

# New column G: "current" header + "NC" for every data row.
$ws.Range("G1").Value = "current"
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 7).Value = "NC"
}

# Match the author's final selection in the saved workbook.
$ws.Range("G31").Select()
